$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4677.857
$ws.Range("I74").Value = 4457.6665
$ws.Range("J74").Value = 5999
$ws.Range("K74").Value = 4457.6665
$ws.Range("L74").Value = 5999
$ws.Range("M74").Value = -3521.6665
$ws.Range("N74").Value = -7871
$ws.Range("H77").Value = 4677.857
$ws.Range("I77").Value = 4457.6665
$ws.Range("J77").Value = 5999
$ws.Range("K77").Value = 22288.3325
$ws.Range("L77").Value = 29995
$ws.Range("M77").Value = -17608.3325
$ws.Range("N77").Value = -39355
$ws.Range("H116").Value = 20992108
$ws.Range("I116").Value = 25761726
$ws.Range("K116").Value = 25761726
$ws.Range("M116").Value = -25758284
$ws.Range("H138").Value = 5547.5923
$ws.Range("I138").Value = 1020.4545
$ws.Range("K138").Value = 3061.3635
$ws.Range("M138").Value = 2078.6365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18079.404
$ws.Range("I32").Value = 17493.977
$ws.Range("K32").Value = 17493.977
$ws.Range("M32").Value = -17206.977
$ws.Range("H61").Value = 2407.4285
$ws.Range("I61").Value = 2407.4285
$ws.Range("K61").Value = 2407.4285
$ws.Range("M61").Value = -2195.4285
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").Value = 2000
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").Value = 10000
$ws.Range("N77").Value = -18736
$ws.Range("H103").Value = 42000
$ws.Range("J103").Value = 42000
$ws.Range("L103").Value = 42000
$ws.Range("N103").Value = -44344
$ws.Range("H136").Value = 2407.4285
$ws.Range("I136").Value = 2407.4285
$ws.Range("K136").Value = 7222.2855
$ws.Range("M136").Value = -4672.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 998
$ws.Range("I37").Value = 998
$ws.Range("K37").Value = 998
$ws.Range("M37").Value = -861
$ws.Range("H86").Value = 3231.5386
$ws.Range("I86").Value = 2523.5557
$ws.Range("K86").Value = 2523.5557
$ws.Range("M86").Value = -1400.5557
$ws.Range("H89").Value = 3231.5386
$ws.Range("I89").Value = 2523.5557
$ws.Range("K89").Value = 12617.7785
$ws.Range("M89").Value = -7001.7785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11990.5
$ws.Range("J99").Value = 4098.3335
$ws.Range("L99").Value = 4098.3335
$ws.Range("N99").Value = -7094.3335
$ws.Range("H126").Value = 11990.5
$ws.Range("J126").Value = 4098.3335
$ws.Range("L126").Value = 12295.0005
$ws.Range("N126").Value = -17235.0005
$ws.Range("H141").Value = 128756.43
$ws.Range("J141").Value = 137560.11
$ws.Range("L141").Value = 137560.11
$ws.Range("N141").Value = -147920.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 126499.5
$ws.Range("I128").Value = 126499.5
$ws.Range("K128").Value = 379498.5
$ws.Range("M128").Value = -374518.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3952.25
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3952.25
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 3952.25
$ws.Range("N97").Value = -4944.25
$ws.Range("H141").Value = 32451.2
$ws.Range("J141").Value = 32451.2
$ws.Range("L141").Value = 32451.2
$ws.Range("N141").Value = -42811.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1999
$ws.Range("J2").Value = 1999
$ws.Range("L2").Value = 1999
$ws.Range("N2").Value = -2223
$ws.Range("H22").Value = 720.7143
$ws.Range("I22").Value = 794.7692
$ws.Range("K22").Value = 794.7692
$ws.Range("M22").Value = -499.7692
$ws.Range("H27").Value = 720.7143
$ws.Range("I27").Value = 794.7692
$ws.Range("K27").Value = 794.7692
$ws.Range("M27").Value = -687.7692
$ws.Range("H40").Value = 5955432.5
$ws.Range("I40").Value = 7354828.5
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 7354828.5
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -7354692.5
$ws.Range("N40").Value = -8272
$ws.Range("H46").Value = 4307.5386
$ws.Range("I46").Value = 2250.5
$ws.Range("K46").Value = 2250.5
$ws.Range("M46").Value = -2062.5
$ws.Range("H68").Value = 5224.4
$ws.Range("I68").Value = 3541.3333
$ws.Range("J68").Value = 7749
$ws.Range("K68").Value = 3541.3333
$ws.Range("L68").Value = 7749
$ws.Range("M68").Value = -2792.3333
$ws.Range("N68").Value = -9247
$ws.Range("H71").Value = 5224.4
$ws.Range("I71").Value = 3541.3333
$ws.Range("J71").Value = 7749
$ws.Range("K71").Value = 17706.6665
$ws.Range("L71").Value = 38745
$ws.Range("M71").Value = -13962.6665
$ws.Range("N71").Value = -46233
$ws.Range("H82").Value = 2073.9546
$ws.Range("I82").Value = 1198.3077
$ws.Range("K82").Value = 1198.3077
$ws.Range("M82").Value = -837.3077000000001
$ws.Range("H85").Value = 2073.9546
$ws.Range("I85").Value = 1198.3077
$ws.Range("K85").Value = 1198.3077
$ws.Range("M85").Value = 49.69229999999993
$ws.Range("H93").Value = 1622
$ws.Range("I93").Value = 1740.6
$ws.Range("J93").Value = 1424.3334
$ws.Range("K93").Value = 1740.6
$ws.Range("L93").Value = 1424.3334
$ws.Range("M93").Value = -492.5999999999999
$ws.Range("N93").Value = -3920.3334
$ws.Range("H122").Value = 8789.200000000001
$ws.Range("I122").Value = 4525.273
$ws.Range("J122").Value = 12139.429
$ws.Range("K122").Value = 13575.819
$ws.Range("L122").Value = 36418.287
$ws.Range("M122").Value = -11125.819
$ws.Range("N122").Value = -41318.287
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6125.75
$ws.Range("J122").Value = 7166.5
$ws.Range("L122").Value = 21499.5
$ws.Range("N122").Value = -26399.5
$ws.Range("H126").Value = 3004
$ws.Range("I126").Value = 3003
$ws.Range("K126").Value = 9009
$ws.Range("M126").Value = -6539
$ws.Range("H138").Value = 98300
$ws.Range("I138").Value = 98300
$ws.Range("K138").Value = 98300
$ws.Range("M138").Value = -93160
$ws.Range("H140").Value = 95272.664
$ws.Range("J140").Value = 95272.664
$ws.Range("L140").Value = 95272.664
$ws.Range("N140").Value = -105632.664
